$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "test"
$ws.Range("B1").Value = "one"

$ws.Range("B1").Select()
